$d = $word.ActiveDocument
$t = $d.Tables.Item(1)

$t.Cell(1, 1).Range.Text  = "84÷3=28, 0"
$t.Cell(1, 2).Range.Text  = "51÷5=10, 1"
$t.Cell(1, 3).Range.Text  = "30÷6=5, 0"
$t.Cell(1, 4).Range.Text  = "35÷3=11, 2"
$t.Cell(1, 5).Range.Text  = "43÷7=6, 1"

$t.Cell(5, 1).Range.Text  = "48÷6=8, 0"
$t.Cell(5, 2).Range.Text  = "71÷7=10, 1"
$t.Cell(5, 3).Range.Text  = "96÷4=24, 0"
$t.Cell(5, 4).Range.Text  = "65÷8=8, 1"
$t.Cell(5, 5).Range.Text  = "90÷6=15, 0"

$t.Cell(9, 1).Range.Text  = "38÷5=7, 3"
$t.Cell(9, 2).Range.Text  = "81÷6=13, 3"
$t.Cell(9, 3).Range.Text  = "72÷5=14, 2"
$t.Cell(9, 4).Range.Text  = "29÷3=9, 2"
$t.Cell(9, 5).Range.Text  = "38÷2=19, 0"

$t.Cell(13, 1).Range.Text = "76÷2=38, 0"
$t.Cell(13, 2).Range.Text = "74÷9=8, 2"
$t.Cell(13, 3).Range.Text = "66÷6=11, 0"
$t.Cell(13, 4).Range.Text = "94÷7=13, 3"
$t.Cell(13, 5).Range.Text = "14÷5=2, 4"

$t.Cell(17, 1).Range.Text = "44÷4=11, 0"
$t.Cell(17, 2).Range.Text = "81÷3=27, 0"
$t.Cell(17, 3).Range.Text = "46÷9=5, 1"
$t.Cell(17, 4).Range.Text = "86÷9=9, 5"
$t.Cell(17, 5).Range.Text = "12÷6=2, 0"
